{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) inside specific resume bullet points by splitting each bullet's\n// single run into multiple runs, bolding + coloring (#2C3E50) the metric\n// substrings while leaving the surrounding text unformatted.\n//\n// Each entry below is one target paragraph's ORIGINAL full text, expressed\n// as an ordered list of [substring, isMetric] pairs whose concatenation\n// reproduces that paragraph's text exactly. Every pair where isMetric is\n// true gets bold + the highlight color; the rest stay plain.\nconst SEGMENTS = [\n  [\n    [\"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from \", false],\n    [\"23%\", true],\n    [\" to \", false],\n    [\"64%\", true],\n  ],\n  [\n    [\"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \", false],\n    [\"\u00b14.2%\", true],\n    [\" to \", false],\n    [\"\u00b12.1%\", true],\n    [\", increasing voter turnout prediction accuracy from \", false],\n    [\"71%\", true],\n    [\" to \", false],\n    [\"87%\", true],\n    [\", and ensuring survey results more closely reflected true population attitudes\", false],\n  ],\n  [\n    [\"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by \", false],\n    [\"73.5%\", true],\n    [\", saving campaigns and organizations \", false],\n    [\"$4.7M\", true],\n    [\" and enabling smaller nonprofits to conduct analysis\", false],\n  ],\n  [\n    [\"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over \", false],\n    [\"$2\", true],\n    [\" trillion\", false],\n  ],\n  [\n    [\"\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by \", false],\n    [\"57%\", true],\n  ],\n  [\n    [\"\u2022 \", false],\n    [\"178%\", true],\n    [\" accuracy improvement in racial classification algorithms\", false],\n  ],\n  [\n    [\"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs \", false],\n    [\"73.5%\", true],\n  ],\n  [\n    [\"\u2022 \", false],\n    [\"$4.7M\", true],\n    [\" savings enabled nonprofit access\", false],\n  ],\n  [\n    [\"\u2022 Platform impact: Built redistricting system serving \", false],\n    [\"12,847\", true],\n    [\" analysts across 89 organizations\", false],\n  ],\n  [\n    [\"\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \", false],\n    [\"\u00b14.2%\", true],\n    [\" to \", false],\n    [\"\u00b12.1%\", true],\n  ],\n  [\n    [\"\u2022 Increased voter turnout prediction accuracy from \", false],\n    [\"71%\", true],\n    [\" to \", false],\n    [\"87%\", true],\n  ],\n];\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Map from a target paragraph's exact original text to its segment list, so\n// we can identify each bullet regardless of its position in the document.\nconst segmentsByText = new Map();\nfor (const segments of SEGMENTS) {\n  const fullText = segments.map(([text]) => text).join(\"\");\n  segmentsByText.set(fullText, segments);\n}\n\nfor (const paragraph of paragraphs.items) {\n  const segments = segmentsByText.get(paragraph.text);\n  if (!segments) continue;\n\n  // Empty the paragraph, then rebuild it run-by-run so each metric\n  // substring lands in its own (bold + colored) run.\n  paragraph.clear();\n  await context.sync();\n\n  for (const [text, isMetric] of segments) {\n    const insertedRange = paragraph.insertText(text, \"End\");\n    if (isMetric) {\n      insertedRange.font.bold = true;\n      insertedRange.font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) inside specific resume bullet points. For each target paragraph we\n# locate the metric substring(s) inside it and bold + color (#2C3E50) just\n# that text, leaving the rest of the paragraph's formatting untouched. Word's\n# COM model stores RGB as 0xBBGGRR in Font.Color, so RgbHexToWdColor() does\n# the byte-order conversion from a normal \"RRGGBB\" hex string.\nfunction RgbHexToWdColor($hex) {\n  $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n  $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n  $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n  return ($b * 65536) + ($g * 256) + $r\n}\n$HighlightColor = RgbHexToWdColor('2C3E50')\n\n$d = $word.ActiveDocument\n\n# Each entry is one target bullet's exact original text plus the ordered list\n# of substrings inside it that must become bold + highlighted.\n$targets = @(\n  @{ Full = '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%'; Bolds = @('23%', '64%') }\n  @{ Full = '\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes'; Bolds = @('\u00b14.2%', '\u00b12.1%', '71%', '87%') }\n  @{ Full = '\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis'; Bolds = @('73.5%', '$4.7M') }\n  @{ Full = '\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion'; Bolds = @('$2') }\n  @{ Full = '\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%'; Bolds = @('57%') }\n  @{ Full = '\u2022 178% accuracy improvement in racial classification algorithms'; Bolds = @('178%') }\n  @{ Full = '\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%'; Bolds = @('73.5%') }\n  @{ Full = '\u2022 $4.7M savings enabled nonprofit access'; Bolds = @('$4.7M') }\n  @{ Full = '\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations'; Bolds = @('12,847') }\n  @{ Full = '\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%'; Bolds = @('\u00b14.2%', '\u00b12.1%') }\n  @{ Full = '\u2022 Increased voter turnout prediction accuracy from 71% to 87%'; Bolds = @('71%', '87%') }\n)\n\nforeach ($target in $targets) {\n  for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $pr = $p.Range()\n\n    # Paragraph Range.Text includes the trailing paragraph mark (CR); trim it\n    # so we can compare against the target's plain text.\n    $ptext = $pr.Text\n    if ($ptext.Length -gt 0 -and $ptext.Substring($ptext.Length - 1) -eq [char]13) {\n      $ptext = $ptext.Substring(0, $ptext.Length - 1)\n    }\n    if ($ptext -ne $target.Full) { continue }\n\n    $base = $pr.Start\n    $cursor = 0\n    foreach ($needle in $target.Bolds) {\n      $relIdx = $ptext.IndexOf($needle, $cursor)\n      if ($relIdx -lt 0) { continue }\n      $s = $base + $relIdx\n      $e = $s + $needle.Length\n      $sub = $d.Range($s, $e)\n      $sub.Font.Bold = 1\n      $sub.Font.Color = $HighlightColor\n      $cursor = $relIdx + $needle.Length\n    }\n    break\n  }\n}\n"}
